$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Requisitos del Proyecto")

# Row 10: "Eliminar juegos de "Mis Juegos"" -> mark as completed, add notes/result
$ws.Range("D10").Value = "Completada"
$ws.Range("G10").Value = 1
$ws.Range("I10").Value = "Los juegos se borran correctamente"
$ws.Range("J10").Value = "OK"

# Row 13: "Recordar contraseña" -> add notes/result
$ws.Range("I13").Value = "Falta la comprobación del cod de verificación pero cumple su objetivo"
$ws.Range("J13").Value = "KO"

# Update selection / view to match final state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L13").Select()
